# Auto-generated edit script for commit 'Add data for 2022-11-11'
# Applies the cell-value updates to the three tiers of sheets:
#   - individual neighborhood sheets
#   - the 'By Neighborhood' roll-up sheet
#   - the 'Citywide Totals' roll-up sheet

$wb = $excel.ActiveWorkbook

# Citywide Totals (sheet1.xml)
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 110
$ws.Range("D3").Value = 126
$ws.Range("B6").Value = 352
$ws.Range("C6").Value = 452
$ws.Range("D6").Value = 383
$ws.Range("E6").Value = 426
$ws.Range("H6").Value = 420
$ws.Range("I6").Value = 471
$ws.Range("B7").Value = 473
$ws.Range("C7").Value = 601
$ws.Range("D7").Value = 601
$ws.Range("E7").Value = 639
$ws.Range("H7").Value = 672
$ws.Range("I7").Value = 786

# By Neighborhood (sheet2.xml)
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("B27").Value = 3
$ws.Range("C32").Value = 38
$ws.Range("I32").Value = 46
$ws.Range("D36").Value = 32
$ws.Range("I41").Value = 9
$ws.Range("E50").Value = 20
$ws.Range("E53").Value = 77
$ws.Range("H61").Value = 6
$ws.Range("I61").Value = 3
$ws.Range("D65").Value = 21
$ws.Range("I72").Value = 10
$ws.Range("I91").Value = 12
$ws.Range("B96").Value = 15
$ws.Range("H97").Value = 3
$ws.Range("B98").Value = 473
$ws.Range("C98").Value = 601
$ws.Range("D98").Value = 601
$ws.Range("E98").Value = 639
$ws.Range("H98").Value = 672
$ws.Range("I98").Value = 786

# Garfield Park (sheet10.xml)
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("C6").Value = 33
$ws.Range("I6").Value = 30
$ws.Range("C7").Value = 38
$ws.Range("I7").Value = 46

# Grand Crossing (sheet11.xml)
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("D3").Value = 8
$ws.Range("D7").Value = 32

# Little Italy, UIC (sheet13.xml)
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("E5").Value = 18
$ws.Range("E6").Value = 20

# Humboldt Park (sheet15.xml)
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 1
$ws.Range("I5").Value = 9

# Woodlawn (sheet21.xml)
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("B5").Value = 10
$ws.Range("B6").Value = 15

# Loop (sheet22.xml)
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("E6").Value = 61
$ws.Range("E7").Value = 77

# West Loop (sheet24.xml)
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I2").Value = 3
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 12

# North Lawndale (sheet31.xml)
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("D5").Value = 20
$ws.Range("D6").Value = 21

# Printers Row (sheet39.xml)
$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("I4").Value = 9
$ws.Range("I5").Value = 10

# Edgewater (sheet44.xml)
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 3

# Wrigleyville (sheet45.xml)
$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 3
